# Update election results for VISEU / CASTRO DAIRE (row 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value  = 156
$ws.Range("I2").Value  = 421
$ws.Range("J2").Value  = 1702
$ws.Range("K2").Value  = 10
$ws.Range("L2").Value  = 436
$ws.Range("M2").Value  = 23
$ws.Range("N2").Value  = 280
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 3
$ws.Range("Q2").Value  = 6
$ws.Range("R2").Value  = 23
$ws.Range("S2").Value  = 194
$ws.Range("T2").Value  = 273
$ws.Range("U2").Value  = 32
$ws.Range("V2").Value  = 2639
$ws.Range("W2").Value  = 1
$ws.Range("X2").Value  = 2618
$ws.Range("Y2").Value  = 3
$ws.Range("Z2").Value  = 43
$ws.Range("AA2").Value = 23
